$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "Your manager has tasked you with creating a Heroku application called forum-node-api using the Heroku CLI. However, due to some shell aliases required for legacy code in the application, they have asked you not to use the heroku create alias to create the application.Which Heroku CLI command should you use to accomplish this?",
        "ques_type": 2,
        "options": [
            "heroku create forum-node-api",
            "heroku git:add forum-node-api",
            "heroku create",
            "heroku app:create forum-node-api"
        ],
        "score": "heroku create forum-node-api"
    },
    {
        "title": "On your local machine, you run your application on port 3000. Now that you are preparing your app to run on Heroku, you want to tell Heroku what port your app will be listening on.How should you communicate this to Heroku?",
        "ques_type": 2,
        "options": [
            "Configure the port in config.md.",
            "Configure the port in Heroku Dashboard.",
            "Set const port = 3000 in Node.js.",
            "Set the PORT environment variable to 3000."
        ],
        "score": "Set const port = 3000 in Node.js."
    },
    {
        "title": "You are trying to deploy your app to Heroku, but it fails during build. You consult the runtime logs, but the error does not appear there.Where should you look to find logs for the unsuccessful build?",
        "ques_type": 2,
        "options": [
            "In the console where you deployed the app.",
            "In the Git log.",
            "In the file /var/log/syslog.",
            "In the Activity tab in the Heroku Dashboard."
        ],
        "score": "In the Activity tab in the Heroku Dashboard."
    },
    {
        "title": "Your Heroku account is protected by a strong password. However, you are worried about the possibility of an attacker leaking your credentials and using them to access your account. How should you harden your account against this kind of threat?",
        "ques_type": 2,
        "options": [
            "Change your password every year.",
            "Salt your password before hashing it.",
            "Use a hash to store your password.",
            "Turn on multi-factor authentication."
        ],
        "score": "Turn on multi-factor authentication."
    }
]
'@

# Remove trailing newline introduced by the here-string
$newText = $newText.TrimEnd("`r", "`n")

# Clear the old row 2 (shared string cell) entirely
$ws.Range("A2").Clear()

# Reset A1's formatting back to default (remove bold font + borders)
$ws.Range("A1").ClearFormats()

# Set A1 to the new, reformatted text
$ws.Range("A1").Value = $newText

# The multi-line text otherwise leaves a stretched custom row height behind;
# autofit brings row 1 back to the sheet's default height.
$ws.Rows.Item(1).AutoFit()
